$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell -> new value map, derived from the "Updated symbol list" diff
# (prices refreshed, Hora 15 -> 16, and rows 41-43 re-ranked with new coins).
$edits = [ordered]@{
    "D2" = "242.82"
    "G2" = "16"
    "D3" = "23.82"
    "G3" = "16"
    "D4" = "5.261"
    "G4" = "16"
    "D5" = "0.05830"
    "G5" = "16"
    "D6" = "6.471"
    "G6" = "16"
    "D7" = "3.344"
    "G7" = "16"
    "D8" = "0.8104"
    "G8" = "16"
    "D9" = "0.8837"
    "G9" = "16"
    "D10" = "0.1384"
    "G10" = "16"
    "D11" = "0.07155"
    "G11" = "16"
    "D12" = "0.03091"
    "G12" = "16"
    "D13" = "0.03062"
    "G13" = "16"
    "D14" = "0.09318"
    "G14" = "16"
    "D15" = "3.831"
    "G15" = "16"
    "D16" = "0.001543"
    "G16" = "16"
    "D17" = "0.04712"
    "G17" = "16"
    "D18" = "0.0006017"
    "G18" = "16"
    "D19" = "0.006169"
    "G19" = "16"
    "G20" = "16"
    "D21" = "0.003847"
    "G21" = "16"
    "G22" = "16"
    "D23" = "3.548"
    "G23" = "16"
    "D24" = "2.166"
    "G24" = "16"
    "D25" = "0.3190"
    "G25" = "16"
    "D26" = "0.1311"
    "G26" = "16"
    "G27" = "16"
    "D28" = "0.0002342"
    "G28" = "16"
    "G29" = "16"
    "G30" = "16"
    "G31" = "16"
    "G32" = "16"
    "G33" = "16"
    "G34" = "16"
    "G35" = "16"
    "G36" = "16"
    "G37" = "16"
    "G38" = "16"
    "G39" = "16"
    "D40" = "0.03787"
    "G40" = "16"
    "B41" = "KickToken"
    "C41" = "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
    "D41" = "0.006355"
    "E41" = "40KickTokenKICK"
    "G41" = "16"
    "B42" = "BKEXToken"
    "C42" = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
    "D42" = "0.1050"
    "E42" = "41BKEXTokenBKK"
    "G42" = "16"
    "B43" = "CEJI"
    "C43" = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
    "D43" = "0.002394"
    "E43" = "42CEJICEJIWorstin24h"
    "G43" = "16"
    "D44" = "0.007033"
    "G44" = "16"
    "D45" = "0.00005342"
    "G45" = "16"
    "D46" = "0.00000000751"
    "G46" = "16"
    "D47" = "0.5405"
    "G47" = "16"
    "D48" = "0.01656"
    "G48" = "16"
    "D49" = "0.00002102"
    "G49" = "16"
    "D50" = "0.0002002"
    "G50" = "16"
    "G51" = "16"
}

foreach ($key in $edits.Keys) {
    $cell = $ws.Range($key)
    $cell.NumberFormat = "@"   # keep values as text, matching the sheet's inline-string cells
    $cell.Value = $edits[$key]
}
